$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace full names with first names only
$ws.Range("A2").Value = "Annaliis"
$ws.Range("A3").Value = "Kertu"
$ws.Range("A4").Value = "Andrea"
$ws.Range("A5").Value = "Tiia"
$ws.Range("A6").Value = "Agnes"
$ws.Range("A7").Value = "Käbi"

# Move the active cell selection from B8 to A8
$ws.Range("A8").Select()
